# Apply the "Change SQLite, XLS and XLSX test data" edit to the
# all-field-types workbook:
#   * drop the "null_column" (column I) entirely
#   * turn the boolean literals in A2/A3 into numeric TRUE()/FALSE()
#     formulas (they keep displaying as TRUE/FALSE via the existing
#     custom number format, but are now numeric 1/0 instead of native
#     booleans)
#   * the values that used to live in the now-removed column I get
#     reused on a new trailing row (row 8), spread across columns B:G
#     (with the "null" value appearing twice)
#   * move the active selection to I9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "null_column" column (header cell I1 plus the I3:I7
# data cells). This shifts nothing else - only column I disappears.
$ws.Columns.Item(9).Delete()

# A2 / A3 were native booleans (t="b"); make them numeric formulas that
# evaluate to the same 1 / 0, while keeping their existing "TRUE/FALSE"
# custom number format (style index 1).
$ws.Range("A2").Formula = "=TRUE()-0"
$ws.Range("A3").Formula = "=FALSE()-0"

# Re-introduce the values that used to be in column I, now on a new
# row 8, shifted over to columns B:G (G8 repeats the "null" string).
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "null"
$ws.Range("D8").Value = "nil"
$ws.Range("E8").Value = "none"
$ws.Range("F8").Value = "n/a"
$ws.Range("G8").Value = "null"

# Match the updated selection recorded in the saved workbook.
$ws.Range("I9").Select() | Out-Null
